# Applies the weekly Choclo price-observation roll described in the commit
# message: a new record is inserted for row 345, every existing record from
# the old row 345 through row 408 shifts down by one row, and the former last
# row (408) lands on a brand-new row 409.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 345: brand-new observation
$ws.Range("D345").Value2 = 44522
$ws.Range("H345").Value2 = 'Dulce o Americano'
$ws.Range("J345").Value2 = 55
$ws.Range("K345").Value2 = 24000
$ws.Range("L345").Value2 = 24000
$ws.Range("M345").Value2 = 24000
$ws.Range("N345").Value2 = '$/malla 70 unidades'
$ws.Range("O345").Value2 = 'Provincia de Colchagua'
$ws.Range("P345").Value2 = 343
$ws.Range("Q345").Value2 = 70

# Rows 346-408: shift each row down into the next (each becomes what used to be the row above it)
$ws.Range("D346").Value2 = 44246
$ws.Range("H346").Value2 = 'Choclero'
$ws.Range("I346").Value2 = 'Primera'
$ws.Range("J346").Value2 = 8300
$ws.Range("K346").Value2 = 200
$ws.Range("L346").Value2 = 230
$ws.Range("M346").Value2 = 214
$ws.Range("N346").Value2 = '$/unidad'
$ws.Range("O346").Value2 = 'Provincia de Talca'
$ws.Range("P346").Value2 = 214
$ws.Range("Q346").Value2 = 1
$ws.Range("D347").Value2 = 44246
$ws.Range("H347").Value2 = 'Dulce o Americano'
$ws.Range("I347").Value2 = 'Primera'
$ws.Range("J347").Value2 = 4500
$ws.Range("K347").Value2 = 150
$ws.Range("L347").Value2 = 150
$ws.Range("M347").Value2 = 150
$ws.Range("N347").Value2 = '$/unidad'
$ws.Range("O347").Value2 = 'Provincia de Talca'
$ws.Range("P347").Value2 = 150
$ws.Range("Q347").Value2 = 1
$ws.Range("D348").Value2 = 44491
$ws.Range("H348").Value2 = 'Dulce o Americano'
$ws.Range("I348").Value2 = 'Primera'
$ws.Range("J348").Value2 = 100
$ws.Range("K348").Value2 = 41000
$ws.Range("L348").Value2 = 42000
$ws.Range("M348").Value2 = 41500
$ws.Range("N348").Value2 = '$/malla 70 unidades'
$ws.Range("O348").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P348").Value2 = 593
$ws.Range("Q348").Value2 = 70
$ws.Range("D349").Value2 = 44272
$ws.Range("H349").Value2 = 'Choclero'
$ws.Range("I349").Value2 = 'Primera'
$ws.Range("J349").Value2 = 5500
$ws.Range("K349").Value2 = 250
$ws.Range("L349").Value2 = 250
$ws.Range("M349").Value2 = 250
$ws.Range("N349").Value2 = '$/unidad'
$ws.Range("O349").Value2 = 'Provincia de Quillota'
$ws.Range("P349").Value2 = 250
$ws.Range("Q349").Value2 = 1
$ws.Range("D350").Value2 = 44272
$ws.Range("H350").Value2 = 'Dulce o Americano'
$ws.Range("I350").Value2 = 'Primera'
$ws.Range("J350").Value2 = 6000
$ws.Range("K350").Value2 = 150
$ws.Range("L350").Value2 = 180
$ws.Range("M350").Value2 = 164
$ws.Range("N350").Value2 = '$/unidad'
$ws.Range("O350").Value2 = 'Provincia de Quillota'
$ws.Range("P350").Value2 = 164
$ws.Range("Q350").Value2 = 1
$ws.Range("D351").Value2 = 44305
$ws.Range("H351").Value2 = 'Choclero'
$ws.Range("I351").Value2 = 'Primera'
$ws.Range("J351").Value2 = 5800
$ws.Range("K351").Value2 = 200
$ws.Range("L351").Value2 = 200
$ws.Range("M351").Value2 = 200
$ws.Range("N351").Value2 = '$/unidad'
$ws.Range("O351").Value2 = 'Provincia de Talca'
$ws.Range("P351").Value2 = 200
$ws.Range("Q351").Value2 = 1
$ws.Range("D352").Value2 = 44305
$ws.Range("H352").Value2 = 'Dulce o Americano'
$ws.Range("I352").Value2 = 'Primera'
$ws.Range("J352").Value2 = 6580
$ws.Range("K352").Value2 = 180
$ws.Range("L352").Value2 = 180
$ws.Range("M352").Value2 = 180
$ws.Range("N352").Value2 = '$/unidad'
$ws.Range("O352").Value2 = 'Provincia de Talca'
$ws.Range("P352").Value2 = 180
$ws.Range("Q352").Value2 = 1
$ws.Range("D353").Value2 = 44166
$ws.Range("H353").Value2 = 'Dulce o Americano'
$ws.Range("I353").Value2 = 'Primera'
$ws.Range("J353").Value2 = 35
$ws.Range("K353").Value2 = 20000
$ws.Range("L353").Value2 = 20000
$ws.Range("M353").Value2 = 20000
$ws.Range("N353").Value2 = '$/malla 50 unidades'
$ws.Range("O353").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P353").Value2 = 400
$ws.Range("Q353").Value2 = 50
$ws.Range("D354").Value2 = 44166
$ws.Range("H354").Value2 = 'Dulce o Americano'
$ws.Range("I354").Value2 = 'Primera'
$ws.Range("J354").Value2 = 1350
$ws.Range("K354").Value2 = 300
$ws.Range("L354").Value2 = 320
$ws.Range("M354").Value2 = 310
$ws.Range("N354").Value2 = '$/unidad'
$ws.Range("O354").Value2 = 'Provincia de Limarí'
$ws.Range("P354").Value2 = 310
$ws.Range("Q354").Value2 = 1
$ws.Range("D355").Value2 = 44225
$ws.Range("H355").Value2 = 'Choclero'
$ws.Range("I355").Value2 = 'Primera'
$ws.Range("J355").Value2 = 4500
$ws.Range("K355").Value2 = 250
$ws.Range("L355").Value2 = 250
$ws.Range("M355").Value2 = 250
$ws.Range("N355").Value2 = '$/unidad'
$ws.Range("O355").Value2 = 'Provincia de Quillota'
$ws.Range("P355").Value2 = 250
$ws.Range("Q355").Value2 = 1
$ws.Range("D356").Value2 = 44225
$ws.Range("H356").Value2 = 'Choclero'
$ws.Range("I356").Value2 = 'Segunda'
$ws.Range("J356").Value2 = 5500
$ws.Range("K356").Value2 = 180
$ws.Range("L356").Value2 = 200
$ws.Range("M356").Value2 = 190
$ws.Range("N356").Value2 = '$/unidad'
$ws.Range("O356").Value2 = 'Provincia de Quillota'
$ws.Range("P356").Value2 = 190
$ws.Range("Q356").Value2 = 1
$ws.Range("D357").Value2 = 44225
$ws.Range("H357").Value2 = 'Dulce o Americano'
$ws.Range("I357").Value2 = 'Primera'
$ws.Range("J357").Value2 = 2400
$ws.Range("K357").Value2 = 200
$ws.Range("L357").Value2 = 200
$ws.Range("M357").Value2 = 200
$ws.Range("N357").Value2 = '$/unidad'
$ws.Range("O357").Value2 = 'Provincia de Quillota'
$ws.Range("P357").Value2 = 200
$ws.Range("Q357").Value2 = 1
$ws.Range("D358").Value2 = 44225
$ws.Range("H358").Value2 = 'Dulce o Americano'
$ws.Range("I358").Value2 = 'Segunda'
$ws.Range("J358").Value2 = 6000
$ws.Range("K358").Value2 = 120
$ws.Range("L358").Value2 = 160
$ws.Range("M358").Value2 = 135
$ws.Range("N358").Value2 = '$/unidad'
$ws.Range("O358").Value2 = 'Provincia de Quillota'
$ws.Range("P358").Value2 = 135
$ws.Range("Q358").Value2 = 1
$ws.Range("D359").Value2 = 44447
$ws.Range("H359").Value2 = 'Dulce o Americano'
$ws.Range("I359").Value2 = 'Primera'
$ws.Range("J359").Value2 = 76
$ws.Range("K359").Value2 = 32000
$ws.Range("L359").Value2 = 33000
$ws.Range("M359").Value2 = 32500
$ws.Range("N359").Value2 = '$/malla 70 unidades'
$ws.Range("O359").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P359").Value2 = 464
$ws.Range("Q359").Value2 = 70
$ws.Range("D360").Value2 = 44425
$ws.Range("H360").Value2 = 'Dulce o Americano'
$ws.Range("I360").Value2 = 'Primera'
$ws.Range("J360").Value2 = 105
$ws.Range("K360").Value2 = 32000
$ws.Range("L360").Value2 = 33000
$ws.Range("M360").Value2 = 32524
$ws.Range("N360").Value2 = '$/malla 70 unidades'
$ws.Range("O360").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P360").Value2 = 465
$ws.Range("Q360").Value2 = 70
$ws.Range("D361").Value2 = 44315
$ws.Range("H361").Value2 = 'Dulce o Americano'
$ws.Range("I361").Value2 = 'Primera'
$ws.Range("J361").Value2 = 4500
$ws.Range("K361").Value2 = 180
$ws.Range("L361").Value2 = 180
$ws.Range("M361").Value2 = 180
$ws.Range("N361").Value2 = '$/unidad'
$ws.Range("O361").Value2 = 'Provincia de Talca'
$ws.Range("P361").Value2 = 180
$ws.Range("Q361").Value2 = 1
$ws.Range("D362").Value2 = 44348
$ws.Range("H362").Value2 = 'Dulce o Americano'
$ws.Range("I362").Value2 = 'Primera'
$ws.Range("J362").Value2 = 4200
$ws.Range("K362").Value2 = 200
$ws.Range("L362").Value2 = 200
$ws.Range("M362").Value2 = 200
$ws.Range("N362").Value2 = '$/unidad'
$ws.Range("O362").Value2 = 'Provincia de Talca'
$ws.Range("P362").Value2 = 200
$ws.Range("Q362").Value2 = 1
$ws.Range("D363").Value2 = 44322
$ws.Range("H363").Value2 = 'Choclero'
$ws.Range("I363").Value2 = 'Primera'
$ws.Range("J363").Value2 = 3400
$ws.Range("K363").Value2 = 230
$ws.Range("L363").Value2 = 250
$ws.Range("M363").Value2 = 239
$ws.Range("N363").Value2 = '$/unidad'
$ws.Range("O363").Value2 = 'Provincia de Quillota'
$ws.Range("P363").Value2 = 239
$ws.Range("Q363").Value2 = 1
$ws.Range("D364").Value2 = 44322
$ws.Range("H364").Value2 = 'Dulce o Americano'
$ws.Range("I364").Value2 = 'Primera'
$ws.Range("J364").Value2 = 4500
$ws.Range("K364").Value2 = 180
$ws.Range("L364").Value2 = 180
$ws.Range("M364").Value2 = 180
$ws.Range("N364").Value2 = '$/unidad'
$ws.Range("O364").Value2 = 'Provincia de Quillota'
$ws.Range("P364").Value2 = 180
$ws.Range("Q364").Value2 = 1
$ws.Range("D365").Value2 = 44495
$ws.Range("H365").Value2 = 'Dulce o Americano'
$ws.Range("I365").Value2 = 'Primera'
$ws.Range("J365").Value2 = 73
$ws.Range("K365").Value2 = 41000
$ws.Range("L365").Value2 = 42000
$ws.Range("M365").Value2 = 41479
$ws.Range("N365").Value2 = '$/malla 70 unidades'
$ws.Range("O365").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P365").Value2 = 593
$ws.Range("Q365").Value2 = 70
$ws.Range("D366").Value2 = 44232
$ws.Range("H366").Value2 = 'Choclero'
$ws.Range("I366").Value2 = 'Primera'
$ws.Range("J366").Value2 = 2500
$ws.Range("K366").Value2 = 250
$ws.Range("L366").Value2 = 250
$ws.Range("M366").Value2 = 250
$ws.Range("N366").Value2 = '$/unidad'
$ws.Range("O366").Value2 = 'Provincia de Quillota'
$ws.Range("P366").Value2 = 250
$ws.Range("Q366").Value2 = 1
$ws.Range("D367").Value2 = 44232
$ws.Range("H367").Value2 = 'Choclero'
$ws.Range("I367").Value2 = 'Segunda'
$ws.Range("J367").Value2 = 850
$ws.Range("K367").Value2 = 180
$ws.Range("L367").Value2 = 180
$ws.Range("M367").Value2 = 180
$ws.Range("N367").Value2 = '$/unidad'
$ws.Range("O367").Value2 = 'Provincia de Quillota'
$ws.Range("P367").Value2 = 180
$ws.Range("Q367").Value2 = 1
$ws.Range("D368").Value2 = 44232
$ws.Range("H368").Value2 = 'Dulce o Americano'
$ws.Range("I368").Value2 = 'Primera'
$ws.Range("J368").Value2 = 4000
$ws.Range("K368").Value2 = 180
$ws.Range("L368").Value2 = 200
$ws.Range("M368").Value2 = 191
$ws.Range("N368").Value2 = '$/unidad'
$ws.Range("O368").Value2 = 'Provincia de Quillota'
$ws.Range("P368").Value2 = 191
$ws.Range("Q368").Value2 = 1
$ws.Range("D369").Value2 = 44232
$ws.Range("H369").Value2 = 'Dulce o Americano'
$ws.Range("I369").Value2 = 'Segunda'
$ws.Range("J369").Value2 = 1600
$ws.Range("K369").Value2 = 140
$ws.Range("L369").Value2 = 140
$ws.Range("M369").Value2 = 140
$ws.Range("N369").Value2 = '$/unidad'
$ws.Range("O369").Value2 = 'Provincia de Quillota'
$ws.Range("P369").Value2 = 140
$ws.Range("Q369").Value2 = 1
$ws.Range("D370").Value2 = 44327
$ws.Range("H370").Value2 = 'Dulce o Americano'
$ws.Range("I370").Value2 = 'Primera'
$ws.Range("J370").Value2 = 3800
$ws.Range("K370").Value2 = 180
$ws.Range("L370").Value2 = 180
$ws.Range("M370").Value2 = 180
$ws.Range("N370").Value2 = '$/unidad'
$ws.Range("O370").Value2 = 'Provincia de Quillota'
$ws.Range("P370").Value2 = 180
$ws.Range("Q370").Value2 = 1
$ws.Range("D371").Value2 = 44510
$ws.Range("H371").Value2 = 'Dulce o Americano'
$ws.Range("I371").Value2 = 'Primera'
$ws.Range("J371").Value2 = 73
$ws.Range("K371").Value2 = 40000
$ws.Range("L371").Value2 = 41000
$ws.Range("M371").Value2 = 40521
$ws.Range("N371").Value2 = '$/malla 70 unidades'
$ws.Range("O371").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P371").Value2 = 579
$ws.Range("Q371").Value2 = 70
$ws.Range("D372").Value2 = 44161
$ws.Range("H372").Value2 = 'Choclero'
$ws.Range("I372").Value2 = 'Primera'
$ws.Range("J372").Value2 = 9500
$ws.Range("K372").Value2 = 400
$ws.Range("L372").Value2 = 400
$ws.Range("M372").Value2 = 400
$ws.Range("N372").Value2 = '$/unidad'
$ws.Range("O372").Value2 = 'Provincia de Limarí'
$ws.Range("P372").Value2 = 400
$ws.Range("Q372").Value2 = 1
$ws.Range("D373").Value2 = 44161
$ws.Range("H373").Value2 = 'Dulce o Americano'
$ws.Range("I373").Value2 = 'Primera'
$ws.Range("J373").Value2 = 85
$ws.Range("K373").Value2 = 19000
$ws.Range("L373").Value2 = 20000
$ws.Range("M373").Value2 = 19471
$ws.Range("N373").Value2 = '$/malla 70 unidades'
$ws.Range("O373").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P373").Value2 = 278
$ws.Range("Q373").Value2 = 70
$ws.Range("D374").Value2 = 44468
$ws.Range("H374").Value2 = 'Dulce o Americano'
$ws.Range("I374").Value2 = 'Primera'
$ws.Range("J374").Value2 = 75
$ws.Range("K374").Value2 = 37000
$ws.Range("L374").Value2 = 38000
$ws.Range("M374").Value2 = 37533
$ws.Range("N374").Value2 = '$/malla 70 unidades'
$ws.Range("O374").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P374").Value2 = 536
$ws.Range("Q374").Value2 = 70
$ws.Range("D375").Value2 = 44517
$ws.Range("H375").Value2 = 'Dulce o Americano'
$ws.Range("I375").Value2 = 'Primera'
$ws.Range("J375").Value2 = 40
$ws.Range("K375").Value2 = 41000
$ws.Range("L375").Value2 = 41000
$ws.Range("M375").Value2 = 41000
$ws.Range("N375").Value2 = '$/malla 70 unidades'
$ws.Range("O375").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P375").Value2 = 586
$ws.Range("Q375").Value2 = 70
$ws.Range("D376").Value2 = 44238
$ws.Range("H376").Value2 = 'Choclero'
$ws.Range("I376").Value2 = 'Primera'
$ws.Range("J376").Value2 = 7400
$ws.Range("K376").Value2 = 170
$ws.Range("L376").Value2 = 300
$ws.Range("M376").Value2 = 242
$ws.Range("N376").Value2 = '$/unidad'
$ws.Range("O376").Value2 = 'Provincia de Quillota'
$ws.Range("P376").Value2 = 242
$ws.Range("Q376").Value2 = 1
$ws.Range("D377").Value2 = 44238
$ws.Range("H377").Value2 = 'Choclero'
$ws.Range("I377").Value2 = 'Segunda'
$ws.Range("J377").Value2 = 2400
$ws.Range("K377").Value2 = 200
$ws.Range("L377").Value2 = 200
$ws.Range("M377").Value2 = 200
$ws.Range("N377").Value2 = '$/unidad'
$ws.Range("O377").Value2 = 'Provincia de Quillota'
$ws.Range("P377").Value2 = 200
$ws.Range("Q377").Value2 = 1
$ws.Range("D378").Value2 = 44238
$ws.Range("H378").Value2 = 'Dulce o Americano'
$ws.Range("I378").Value2 = 'Primera'
$ws.Range("J378").Value2 = 2800
$ws.Range("K378").Value2 = 200
$ws.Range("L378").Value2 = 200
$ws.Range("M378").Value2 = 200
$ws.Range("N378").Value2 = '$/unidad'
$ws.Range("O378").Value2 = 'Provincia de Quillota'
$ws.Range("P378").Value2 = 200
$ws.Range("Q378").Value2 = 1
$ws.Range("D379").Value2 = 44238
$ws.Range("H379").Value2 = 'Dulce o Americano'
$ws.Range("I379").Value2 = 'Segunda'
$ws.Range("J379").Value2 = 1500
$ws.Range("K379").Value2 = 150
$ws.Range("L379").Value2 = 150
$ws.Range("M379").Value2 = 150
$ws.Range("N379").Value2 = '$/unidad'
$ws.Range("O379").Value2 = 'Provincia de Quillota'
$ws.Range("P379").Value2 = 150
$ws.Range("Q379").Value2 = 1
$ws.Range("D380").Value2 = 44391
$ws.Range("H380").Value2 = 'Dulce o Americano'
$ws.Range("I380").Value2 = 'Primera'
$ws.Range("J380").Value2 = 65
$ws.Range("K380").Value2 = 24000
$ws.Range("L380").Value2 = 24000
$ws.Range("M380").Value2 = 24000
$ws.Range("N380").Value2 = '$/malla 70 unidades'
$ws.Range("O380").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P380").Value2 = 343
$ws.Range("Q380").Value2 = 70
$ws.Range("D381").Value2 = 44236
$ws.Range("H381").Value2 = 'Choclero'
$ws.Range("I381").Value2 = 'Primera'
$ws.Range("J381").Value2 = 7000
$ws.Range("K381").Value2 = 250
$ws.Range("L381").Value2 = 280
$ws.Range("M381").Value2 = 265
$ws.Range("N381").Value2 = '$/unidad'
$ws.Range("O381").Value2 = 'Provincia de Talca'
$ws.Range("P381").Value2 = 265
$ws.Range("Q381").Value2 = 1
$ws.Range("D382").Value2 = 44236
$ws.Range("H382").Value2 = 'Choclero'
$ws.Range("I382").Value2 = 'Segunda'
$ws.Range("J382").Value2 = 3600
$ws.Range("K382").Value2 = 200
$ws.Range("L382").Value2 = 200
$ws.Range("M382").Value2 = 200
$ws.Range("N382").Value2 = '$/unidad'
$ws.Range("O382").Value2 = 'Provincia de Talca'
$ws.Range("P382").Value2 = 200
$ws.Range("Q382").Value2 = 1
$ws.Range("D383").Value2 = 44236
$ws.Range("H383").Value2 = 'Dulce o Americano'
$ws.Range("I383").Value2 = 'Primera'
$ws.Range("J383").Value2 = 8300
$ws.Range("K383").Value2 = 180
$ws.Range("L383").Value2 = 200
$ws.Range("M383").Value2 = 190
$ws.Range("N383").Value2 = '$/unidad'
$ws.Range("O383").Value2 = 'Provincia de Talca'
$ws.Range("P383").Value2 = 190
$ws.Range("Q383").Value2 = 1
$ws.Range("D384").Value2 = 44389
$ws.Range("H384").Value2 = 'Dulce o Americano'
$ws.Range("I384").Value2 = 'Primera'
$ws.Range("J384").Value2 = 60
$ws.Range("K384").Value2 = 25000
$ws.Range("L384").Value2 = 25000
$ws.Range("M384").Value2 = 25000
$ws.Range("N384").Value2 = '$/malla 70 unidades'
$ws.Range("O384").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P384").Value2 = 357
$ws.Range("Q384").Value2 = 70
$ws.Range("D385").Value2 = 44251
$ws.Range("H385").Value2 = 'Choclero'
$ws.Range("I385").Value2 = 'Primera'
$ws.Range("J385").Value2 = 5200
$ws.Range("K385").Value2 = 280
$ws.Range("L385").Value2 = 300
$ws.Range("M385").Value2 = 290
$ws.Range("N385").Value2 = '$/unidad'
$ws.Range("O385").Value2 = 'Provincia de Quillota'
$ws.Range("P385").Value2 = 290
$ws.Range("Q385").Value2 = 1
$ws.Range("D386").Value2 = 44251
$ws.Range("H386").Value2 = 'Choclero'
$ws.Range("I386").Value2 = 'Segunda'
$ws.Range("J386").Value2 = 2800
$ws.Range("K386").Value2 = 200
$ws.Range("L386").Value2 = 200
$ws.Range("M386").Value2 = 200
$ws.Range("N386").Value2 = '$/unidad'
$ws.Range("O386").Value2 = 'Provincia de Quillota'
$ws.Range("P386").Value2 = 200
$ws.Range("Q386").Value2 = 1
$ws.Range("D387").Value2 = 44251
$ws.Range("H387").Value2 = 'Dulce o Americano'
$ws.Range("I387").Value2 = 'Primera'
$ws.Range("J387").Value2 = 3500
$ws.Range("K387").Value2 = 180
$ws.Range("L387").Value2 = 180
$ws.Range("M387").Value2 = 180
$ws.Range("N387").Value2 = '$/unidad'
$ws.Range("O387").Value2 = 'Provincia de Quillota'
$ws.Range("P387").Value2 = 180
$ws.Range("Q387").Value2 = 1
$ws.Range("D388").Value2 = 44340
$ws.Range("H388").Value2 = 'Choclero'
$ws.Range("I388").Value2 = 'Primera'
$ws.Range("J388").Value2 = 3600
$ws.Range("K388").Value2 = 300
$ws.Range("L388").Value2 = 300
$ws.Range("M388").Value2 = 300
$ws.Range("N388").Value2 = '$/unidad'
$ws.Range("O388").Value2 = 'Provincia de Quillota'
$ws.Range("P388").Value2 = 300
$ws.Range("Q388").Value2 = 1
$ws.Range("D389").Value2 = 44340
$ws.Range("H389").Value2 = 'Dulce o Americano'
$ws.Range("I389").Value2 = 'Primera'
$ws.Range("J389").Value2 = 2200
$ws.Range("K389").Value2 = 250
$ws.Range("L389").Value2 = 250
$ws.Range("M389").Value2 = 250
$ws.Range("N389").Value2 = '$/unidad'
$ws.Range("O389").Value2 = 'Provincia de Quillota'
$ws.Range("P389").Value2 = 250
$ws.Range("Q389").Value2 = 1
$ws.Range("D390").Value2 = 44340
$ws.Range("H390").Value2 = 'Dulce o Americano'
$ws.Range("I390").Value2 = 'Segunda'
$ws.Range("J390").Value2 = 1800
$ws.Range("K390").Value2 = 200
$ws.Range("L390").Value2 = 200
$ws.Range("M390").Value2 = 200
$ws.Range("N390").Value2 = '$/unidad'
$ws.Range("O390").Value2 = 'Provincia de Quillota'
$ws.Range("P390").Value2 = 200
$ws.Range("Q390").Value2 = 1
$ws.Range("D391").Value2 = 44515
$ws.Range("H391").Value2 = 'Dulce o Americano'
$ws.Range("I391").Value2 = 'Primera'
$ws.Range("J391").Value2 = 73
$ws.Range("K391").Value2 = 40000
$ws.Range("L391").Value2 = 41000
$ws.Range("M391").Value2 = 40521
$ws.Range("N391").Value2 = '$/malla 70 unidades'
$ws.Range("O391").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P391").Value2 = 579
$ws.Range("Q391").Value2 = 70
$ws.Range("D392").Value2 = 44330
$ws.Range("H392").Value2 = 'Dulce o Americano'
$ws.Range("I392").Value2 = 'Primera'
$ws.Range("J392").Value2 = 3500
$ws.Range("K392").Value2 = 180
$ws.Range("L392").Value2 = 180
$ws.Range("M392").Value2 = 180
$ws.Range("N392").Value2 = '$/unidad'
$ws.Range("O392").Value2 = 'Provincia de Quillota'
$ws.Range("P392").Value2 = 180
$ws.Range("Q392").Value2 = 1
$ws.Range("D393").Value2 = 44432
$ws.Range("H393").Value2 = 'Dulce o Americano'
$ws.Range("I393").Value2 = 'Primera'
$ws.Range("J393").Value2 = 38
$ws.Range("K393").Value2 = 32000
$ws.Range("L393").Value2 = 32000
$ws.Range("M393").Value2 = 32000
$ws.Range("N393").Value2 = '$/malla 70 unidades'
$ws.Range("O393").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P393").Value2 = 457
$ws.Range("Q393").Value2 = 70
$ws.Range("D394").Value2 = 44181
$ws.Range("H394").Value2 = 'Dulce o Americano'
$ws.Range("I394").Value2 = 'Primera'
$ws.Range("J394").Value2 = 2800
$ws.Range("K394").Value2 = 300
$ws.Range("L394").Value2 = 300
$ws.Range("M394").Value2 = 300
$ws.Range("N394").Value2 = '$/unidad'
$ws.Range("O394").Value2 = 'Provincia de Talca'
$ws.Range("P394").Value2 = 300
$ws.Range("Q394").Value2 = 1
$ws.Range("D395").Value2 = 44181
$ws.Range("H395").Value2 = 'Dulce o Americano'
$ws.Range("I395").Value2 = 'Segunda'
$ws.Range("J395").Value2 = 2600
$ws.Range("K395").Value2 = 200
$ws.Range("L395").Value2 = 200
$ws.Range("M395").Value2 = 200
$ws.Range("N395").Value2 = '$/unidad'
$ws.Range("O395").Value2 = 'Provincia de Talca'
$ws.Range("P395").Value2 = 200
$ws.Range("Q395").Value2 = 1
$ws.Range("D396").Value2 = 44194
$ws.Range("H396").Value2 = 'Choclero'
$ws.Range("I396").Value2 = 'Primera'
$ws.Range("J396").Value2 = 5300
$ws.Range("K396").Value2 = 350
$ws.Range("L396").Value2 = 400
$ws.Range("M396").Value2 = 383
$ws.Range("N396").Value2 = '$/unidad'
$ws.Range("O396").Value2 = 'Provincia de Talca'
$ws.Range("P396").Value2 = 383
$ws.Range("Q396").Value2 = 1
$ws.Range("D397").Value2 = 44194
$ws.Range("H397").Value2 = 'Choclero'
$ws.Range("I397").Value2 = 'Segunda'
$ws.Range("J397").Value2 = 2100
$ws.Range("K397").Value2 = 250
$ws.Range("L397").Value2 = 250
$ws.Range("M397").Value2 = 250
$ws.Range("N397").Value2 = '$/unidad'
$ws.Range("O397").Value2 = 'Provincia de Talca'
$ws.Range("P397").Value2 = 250
$ws.Range("Q397").Value2 = 1
$ws.Range("D398").Value2 = 44194
$ws.Range("H398").Value2 = 'Dulce o Americano'
$ws.Range("I398").Value2 = 'Primera'
$ws.Range("J398").Value2 = 3800
$ws.Range("K398").Value2 = 300
$ws.Range("L398").Value2 = 300
$ws.Range("M398").Value2 = 300
$ws.Range("N398").Value2 = '$/unidad'
$ws.Range("O398").Value2 = 'Provincia de Talca'
$ws.Range("P398").Value2 = 300
$ws.Range("Q398").Value2 = 1
$ws.Range("D399").Value2 = 44194
$ws.Range("H399").Value2 = 'Dulce o Americano'
$ws.Range("I399").Value2 = 'Segunda'
$ws.Range("J399").Value2 = 1100
$ws.Range("K399").Value2 = 200
$ws.Range("L399").Value2 = 200
$ws.Range("M399").Value2 = 200
$ws.Range("N399").Value2 = '$/unidad'
$ws.Range("O399").Value2 = 'Provincia de Talca'
$ws.Range("P399").Value2 = 200
$ws.Range("Q399").Value2 = 1
$ws.Range("D400").Value2 = 44271
$ws.Range("H400").Value2 = 'Choclero'
$ws.Range("I400").Value2 = 'Primera'
$ws.Range("J400").Value2 = 11300
$ws.Range("K400").Value2 = 200
$ws.Range("L400").Value2 = 230
$ws.Range("M400").Value2 = 215
$ws.Range("N400").Value2 = '$/unidad'
$ws.Range("O400").Value2 = 'Provincia de Quillota'
$ws.Range("P400").Value2 = 215
$ws.Range("Q400").Value2 = 1
$ws.Range("D401").Value2 = 44271
$ws.Range("H401").Value2 = 'Dulce o Americano'
$ws.Range("I401").Value2 = 'Primera'
$ws.Range("J401").Value2 = 4500
$ws.Range("K401").Value2 = 150
$ws.Range("L401").Value2 = 150
$ws.Range("M401").Value2 = 150
$ws.Range("N401").Value2 = '$/unidad'
$ws.Range("O401").Value2 = 'Provincia de Quillota'
$ws.Range("P401").Value2 = 150
$ws.Range("Q401").Value2 = 1
$ws.Range("D402").Value2 = 44307
$ws.Range("H402").Value2 = 'Dulce o Americano'
$ws.Range("I402").Value2 = 'Primera'
$ws.Range("J402").Value2 = 3500
$ws.Range("K402").Value2 = 180
$ws.Range("L402").Value2 = 180
$ws.Range("M402").Value2 = 180
$ws.Range("N402").Value2 = '$/unidad'
$ws.Range("O402").Value2 = 'Provincia de Talca'
$ws.Range("P402").Value2 = 180
$ws.Range("Q402").Value2 = 1
$ws.Range("D403").Value2 = 44400
$ws.Range("H403").Value2 = 'Dulce o Americano'
$ws.Range("I403").Value2 = 'Primera'
$ws.Range("J403").Value2 = 50
$ws.Range("K403").Value2 = 24000
$ws.Range("L403").Value2 = 24000
$ws.Range("M403").Value2 = 24000
$ws.Range("N403").Value2 = '$/malla 70 unidades'
$ws.Range("O403").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P403").Value2 = 343
$ws.Range("Q403").Value2 = 70
$ws.Range("D404").Value2 = 44309
$ws.Range("H404").Value2 = 'Choclero'
$ws.Range("I404").Value2 = 'Primera'
$ws.Range("J404").Value2 = 3500
$ws.Range("K404").Value2 = 200
$ws.Range("L404").Value2 = 200
$ws.Range("M404").Value2 = 200
$ws.Range("N404").Value2 = '$/unidad'
$ws.Range("O404").Value2 = 'Provincia de Talca'
$ws.Range("P404").Value2 = 200
$ws.Range("Q404").Value2 = 1
$ws.Range("D405").Value2 = 44508
$ws.Range("H405").Value2 = 'Dulce o Americano'
$ws.Range("I405").Value2 = 'Primera'
$ws.Range("J405").Value2 = 85
$ws.Range("K405").Value2 = 40000
$ws.Range("L405").Value2 = 41000
$ws.Range("M405").Value2 = 40529
$ws.Range("N405").Value2 = '$/malla 70 unidades'
$ws.Range("O405").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P405").Value2 = 579
$ws.Range("Q405").Value2 = 70
$ws.Range("D406").Value2 = 44201
$ws.Range("H406").Value2 = 'Choclero'
$ws.Range("I406").Value2 = 'Primera'
$ws.Range("J406").Value2 = 3300
$ws.Range("K406").Value2 = 330
$ws.Range("L406").Value2 = 330
$ws.Range("M406").Value2 = 330
$ws.Range("N406").Value2 = '$/unidad'
$ws.Range("O406").Value2 = 'Provincia de Talca'
$ws.Range("P406").Value2 = 330
$ws.Range("Q406").Value2 = 1
$ws.Range("D407").Value2 = 44201
$ws.Range("H407").Value2 = 'Choclero'
$ws.Range("I407").Value2 = 'Segunda'
$ws.Range("J407").Value2 = 3200
$ws.Range("K407").Value2 = 250
$ws.Range("L407").Value2 = 250
$ws.Range("M407").Value2 = 250
$ws.Range("N407").Value2 = '$/unidad'
$ws.Range("O407").Value2 = 'Provincia de Talca'
$ws.Range("P407").Value2 = 250
$ws.Range("Q407").Value2 = 1
$ws.Range("D408").Value2 = 44201
$ws.Range("H408").Value2 = 'Dulce o Americano'
$ws.Range("I408").Value2 = 'Primera'
$ws.Range("J408").Value2 = 3950
$ws.Range("K408").Value2 = 300
$ws.Range("L408").Value2 = 300
$ws.Range("M408").Value2 = 300
$ws.Range("N408").Value2 = '$/unidad'
$ws.Range("O408").Value2 = 'Provincia de Talca'
$ws.Range("P408").Value2 = 300
$ws.Range("Q408").Value2 = 1

# Row 409: new row appended, containing what used to be row 408's observation
$ws.Range("A409").Value2 = 3
$ws.Range("B409").Value2 = 'Femacal de La Calera'
$ws.Range("C409").Value2 = 'Coquimbo'
$ws.Range("D409").Value2 = 44201
$ws.Range("E409").Value2 = 5
$ws.Range("F409").Value2 = 100112024
$ws.Range("G409").Value2 = 'Choclo'
$ws.Range("H409").Value2 = 'Dulce o Americano'
$ws.Range("I409").Value2 = 'Segunda'
$ws.Range("J409").Value2 = 3500
$ws.Range("K409").Value2 = 200
$ws.Range("L409").Value2 = 200
$ws.Range("M409").Value2 = 200
$ws.Range("N409").Value2 = '$/unidad'
$ws.Range("O409").Value2 = 'Provincia de Talca'
$ws.Range("P409").Value2 = 200
$ws.Range("Q409").Value2 = 1
$ws.Range("R409").Value2 = 'Hortaliza'

# Match the date-time number format used by column D throughout the table
$ws.Range("D409").NumberFormat = $ws.Range("D408").NumberFormat
